# Update the Premier League Russia standings table (rows 2-17, columns B-G).
# Column A (#) stays the same; B is team name, C=Cartoes, D=Escanteios,
# E=1.5+, F=2.5+, G=Med. Gols. All values in this sheet are plain text
# (even the numeric-looking ones), so we force Text number format before
# writing the values, then clear the (no longer needed) number format back
# off the cells so they keep the sheet's original "no explicit style" look.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("Zenit",              "1.3", "6.4",  "83%",  "50%", "2.67"),
  @("Krasnodar",          "2.3", "5.1",  "83%",  "38%", "2.54"),
  @("Dinamo Moskva",      "1.8", "5.8",  "84%",  "60%", "3.00"),
  @("Lokomotiv Moskva",   "2.4", "6.0",  "84%",  "52%", "3.04"),
  @("Spartak Moskva",     "2.1", "5.2",  "75%",  "50%", "2.63"),
  @("CSKA Moskva",        "1.8", "4.5",  "83%",  "58%", "3.21"),
  @("Krylya Sovetov",     "2.1", "4.8",  "100%", "67%", "3.21"),
  @("Rostov",             "2.1", "5.4",  "79%",  "62%", "3.08"),
  @("Rubin Kazan",        "1.7", "4.6",  "63%",  "38%", "2.13"),
  @("FK Nizjni Novgorod", "2.0", "3.6",  "60%",  "28%", "2.28"),
  @("Orenburg",           "2.5", "5.4",  "84%",  "40%", "2.44"),
  @("Fakel",              "2.1", "4.7",  "66%",  "28%", "1.92"),
  @("Ural",               "2.0", "5.6",  "67%",  "54%", "2.58"),
  @("Baltika",            "1.5", "4.6",  "54%",  "33%", "2.13"),
  @("Akhmat Grozny",      "2.5", "4.3",  "75%",  "50%", "2.50"),
  @("Sochi",              "2.1", "4.1",  "79%",  "50%", "2.71")
)

$rowStart = 2
$rowEnd = $rowStart + $data.Length - 1
$rng = $ws.Range("B" + $rowStart + ":G" + $rowEnd)

# Force every cell in the block to Text format first so values such as
# "83%" or "1.3" are kept as literal strings instead of being parsed into
# numbers/percentages by Excel's input parser.
$rng.NumberFormat = "@"

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $rowStart + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
}

# Drop the temporary Text number format again so the cells end up with no
# explicit style, matching the rest of the data rows in the sheet.
$rng.ClearFormats()
